$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.675.08"
$ws.Range("E2").Value = "  +0.36%  "
$ws.Range("D3").Value = "2.307.91"
$ws.Range("E3").Value = "  +1.12%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "301.61"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.72%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.76"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.52%  "
$ws.Range("E7").Value = "  +1.11%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.491"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.21%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.42"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.52%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.05"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.84%  "
$ws.Range("E12").Value = "  +0.64%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.119"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.55%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.70"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.43%  "
$ws.Range("D15").Value = "2.668.11"
$ws.Range("E15").Value = "  +1.21%  "
$ws.Range("D16").Value = "2.302.37"
$ws.Range("E16").Value = "  +1.08%  "
$ws.Range("E17").Value = "  +1.50%  "
$ws.Range("D18").Value = "42.616.27"
$ws.Range("E18").Value = "  +0.46%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.20"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.19%  "
$ws.Range("E20").Value = "  +1.62%  "
$ws.Range("D21").Value = "0.0₃0888"
$ws.Range("E21").Value = "  +0.12%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.99"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.38%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.27"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.25%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "235.49"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.10%  "
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.41"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.30%  "
$ws.Range("E27").Value = "  -1.58%  "
$ws.Range("E28").Value = "  -1.04%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "165.49"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.17%  "
$ws.Range("E30").Value = "  +1.94%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "32.25"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.04%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.99"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.37%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "17.55"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.26%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.44"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.74%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0698"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.54%  "
$ws.Range("E37").Value = "  -0.86%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0996"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.48%  "
$ws.Range("E39").Value = "  +1.12%  "
$ws.Range("E40").Value = "  -0.12%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.68"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.43%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "19.81"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +11.28%  "
$ws.Range("D43").Value = "1.951.88"
$ws.Range("E43").Value = "  -1.85%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.43"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.13%  "
$ws.Range("E45").Value = "  +1.19%  "
$ws.Range("E46").Value = "  +2.62%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.72"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.64%  "
$ws.Range("D48").Value = "2.536.86"
$ws.Range("E48").Value = "  +1.29%  "
$ws.Range("B49").Value = "HuobiToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.79"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.64%  "
$ws.Range("B50").Value = "MultiversX"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "53.22"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.33%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "71.99"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.15%  "
